# Apply corrected stock quantities/values (and recomputed sub/grand totals)
# to the CryCompanywiseStockReport worksheet, matching the target revision.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(71, 6).Value = 315
$ws.Cells.Item(71, 7).Value = 20065.5
$ws.Cells.Item(73, 6).Value = 68
$ws.Cells.Item(73, 7).Value = 5370.64
$ws.Cells.Item(74, 6).Value = 141
$ws.Cells.Item(74, 7).Value = 19789.35
$ws.Cells.Item(77, 6).Value = 245
$ws.Cells.Item(77, 7).Value = 11451.3
$ws.Cells.Item(86, 6).Value = 54
$ws.Cells.Item(86, 7).Value = 6775.38
$ws.Cells.Item(90, 2).Value = 172935.59
$ws.Cells.Item(102, 6).Value = 4
$ws.Cells.Item(102, 7).Value = 197.92
$ws.Cells.Item(104, 2).Value = 165.9
$ws.Cells.Item(115, 6).Value = 191
$ws.Cells.Item(115, 7).Value = 18490.71
$ws.Cells.Item(117, 2).Value = 12348.59
$ws.Cells.Item(135, 6).Value = 21
$ws.Cells.Item(135, 7).Value = 651.63
$ws.Cells.Item(138, 2).Value = 2482.77
$ws.Cells.Item(144, 6).Value = 1012
$ws.Cells.Item(144, 7).Value = 8551.4
$ws.Cells.Item(145, 6).Value = 412
$ws.Cells.Item(145, 7).Value = 3291.88
$ws.Cells.Item(146, 6).Value = 18
$ws.Cells.Item(146, 7).Value = 1515.42
$ws.Cells.Item(147, 2).Value = 13358.7
$ws.Cells.Item(149, 6).Value = 224
$ws.Cells.Item(149, 7).Value = 14515.2
$ws.Cells.Item(150, 6).Value = 32
$ws.Cells.Item(150, 7).Value = 1487.68
$ws.Cells.Item(156, 2).Value = 30834.17
$ws.Cells.Item(192, 2).Value = 48706
$ws.Cells.Item(192, 5).Value = 39.8
$ws.Cells.Item(192, 6).Value = -144
$ws.Cells.Item(192, 7).Value = -4795.2
$ws.Cells.Item(193, 2).Value = 64973
$ws.Cells.Item(193, 5).Value = 35.4
$ws.Cells.Item(193, 6).Value = 2
$ws.Cells.Item(193, 7).Value = 66.59999999999999
$ws.Cells.Item(203, 6).Value = 54
$ws.Cells.Item(203, 7).Value = 1088.64
$ws.Cells.Item(216, 2).Value = 37552.52
$ws.Cells.Item(229, 6).Value = 54
$ws.Cells.Item(229, 7).Value = 7747.92
$ws.Cells.Item(255, 6).Value = 534
$ws.Cells.Item(255, 7).Value = 91490.22
$ws.Cells.Item(260, 2).Value = 177492.1
$ws.Cells.Item(280, 6).Value = 131
$ws.Cells.Item(280, 7).Value = 22157.34
$ws.Cells.Item(291, 6).Value = 108
$ws.Cells.Item(291, 7).Value = 4645.08
$ws.Cells.Item(293, 6).Value = 32
$ws.Cells.Item(293, 7).Value = 2250.24
$ws.Cells.Item(296, 6).Value = 53
$ws.Cells.Item(296, 7).Value = 1123.6
$ws.Cells.Item(304, 2).Value = 170342.58
$ws.Cells.Item(306, 6).Value = 61
$ws.Cells.Item(306, 7).Value = 1287.71
$ws.Cells.Item(309, 2).Value = 1709.48
$ws.Cells.Item(320, 6).Value = 41
$ws.Cells.Item(320, 7).Value = 2814.65
$ws.Cells.Item(322, 2).Value = 58047
$ws.Cells.Item(322, 4).Value = 105.54
$ws.Cells.Item(322, 5).Value = 126.1
$ws.Cells.Item(322, 6).Value = 39
$ws.Cells.Item(322, 7).Value = 4116.06
$ws.Cells.Item(323, 2).Value = 47097
$ws.Cells.Item(323, 4).Value = 112.28
$ws.Cells.Item(323, 5).Value = 134.16
$ws.Cells.Item(323, 6).Value = 15
$ws.Cells.Item(323, 7).Value = 1684.2
$ws.Cells.Item(326, 6).Value = 62
$ws.Cells.Item(326, 7).Value = 1843.88
$ws.Cells.Item(330, 2).Value = 26815.27
$ws.Cells.Item(338, 6).Value = 75
$ws.Cells.Item(338, 7).Value = 1777.5
$ws.Cells.Item(343, 6).Value = 35
$ws.Cells.Item(343, 7).Value = 2518.95
$ws.Cells.Item(345, 6).Value = 43
$ws.Cells.Item(345, 7).Value = 2640.63
$ws.Cells.Item(346, 2).Value = 25003.71
$ws.Cells.Item(354, 6).Value = 12
$ws.Cells.Item(354, 7).Value = 823.08
$ws.Cells.Item(358, 2).Value = 34846.07
$ws.Cells.Item(442, 2).Value = 64810
$ws.Cells.Item(442, 5).Value = 291.22
$ws.Cells.Item(442, 6).Value = 4
$ws.Cells.Item(442, 7).Value = 1095.68
$ws.Cells.Item(443, 2).Value = 53319
$ws.Cells.Item(443, 5).Value = 310.64
$ws.Cells.Item(443, 6).Value = -6
$ws.Cells.Item(443, 7).Value = -1643.52
$ws.Cells.Item(455, 6).Value = 43
$ws.Cells.Item(455, 7).Value = 2735.23
$ws.Cells.Item(460, 2).Value = 12661.57
$ws.Cells.Item(463, 2).Value = 64833
$ws.Cells.Item(463, 5).Value = 34.9
$ws.Cells.Item(463, 6).Value = 95
$ws.Cells.Item(463, 7).Value = 3118.85
$ws.Cells.Item(464, 2).Value = 60025
$ws.Cells.Item(464, 5).Value = 37.22
$ws.Cells.Item(464, 6).Value = -98
$ws.Cells.Item(464, 7).Value = -3217.34
$ws.Cells.Item(482, 6).Value = 41
$ws.Cells.Item(482, 7).Value = 2430.07
$ws.Cells.Item(485, 6).Value = 10
$ws.Cells.Item(485, 7).Value = 1754.7
$ws.Cells.Item(488, 2).Value = 29161.39
$ws.Cells.Item(542, 6).Value = 46
$ws.Cells.Item(542, 7).Value = 5958.38
$ws.Cells.Item(547, 2).Value = 16948.5
$ws.Cells.Item(556, 6).Value = 1
$ws.Cells.Item(556, 7).Value = 114.86
$ws.Cells.Item(560, 2).Value = 4195.03
$ws.Cells.Item(572, 2).Value = 65079
$ws.Cells.Item(572, 6).Value = 6
$ws.Cells.Item(572, 7).Value = 245.22
$ws.Cells.Item(573, 2).Value = 65362
$ws.Cells.Item(573, 6).Value = 20
$ws.Cells.Item(573, 7).Value = 817.4
$ws.Cells.Item(578, 6).Value = 67
$ws.Cells.Item(578, 7).Value = 3342.63
$ws.Cells.Item(583, 2).Value = 14022.29
$ws.Cells.Item(599, 6).Value = 1466
$ws.Cells.Item(599, 7).Value = 239119.26
$ws.Cells.Item(602, 6).Value = 321
$ws.Cells.Item(602, 7).Value = 46432.65
$ws.Cells.Item(606, 2).Value = 394456.3
$ws.Cells.Item(613, 6).Value = 133
$ws.Cells.Item(613, 7).Value = 21168.28
$ws.Cells.Item(618, 2).Value = 42930.25
$ws.Cells.Item(619, 2).Value = 1657907.61
$ws.Cells.Item(620, 2).Value = 1657907.61
